$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Deuxième Test avec commit." paragraph: was plain red; becomes bold,
#    dark-gray theme color (Text 1, Lighter 15%) and 16pt (sz=32 half-points).
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(3)
$rng2 = $p2.Range
$xml2 = @"
<?xml version='1.0'?>
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>
<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>
<pkg:xmlData>
<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:body>
<w:p>
<w:pPr><w:rPr><w:b/><w:color w:val='262626' w:themeColor='text1' w:themeTint='D9'/><w:sz w:val='32'/></w:rPr></w:pPr>
<w:r><w:rPr><w:b/><w:color w:val='262626' w:themeColor='text1' w:themeTint='D9'/><w:sz w:val='32'/></w:rPr><w:t>Deuxième Test avec commit.</w:t></w:r>
<w:r><w:rPr><w:b/><w:color w:val='262626' w:themeColor='text1' w:themeTint='D9'/><w:sz w:val='32'/></w:rPr><w:tab/></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$rng2.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 2) "Troisième Test avec commit." paragraph: already bold+red; becomes bold
#    with the same dark-gray theme color and 16pt.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(4)
$rng3 = $p3.Range
$xml3 = @"
<?xml version='1.0'?>
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>
<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>
<pkg:xmlData>
<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:body>
<w:p>
<w:pPr><w:rPr><w:b/><w:color w:val='262626' w:themeColor='text1' w:themeTint='D9'/><w:sz w:val='32'/></w:rPr></w:pPr>
<w:r><w:rPr><w:b/><w:color w:val='262626' w:themeColor='text1' w:themeTint='D9'/><w:sz w:val='32'/></w:rPr><w:t>Troisième</w:t></w:r>
<w:r><w:rPr><w:b/><w:color w:val='262626' w:themeColor='text1' w:themeTint='D9'/><w:sz w:val='32'/></w:rPr><w:t xml:space='preserve'> Test avec commit.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$rng3.InsertXML($xml3)

# ---------------------------------------------------------------------------
# 3) Move the (hidden) "_GoBack" bookmark from the trailing empty paragraph
#    up into the "CONCLUSION GENERALE ET PERSPECTIVES" title, splitting that
#    run right before "PERSPECTIVES". Word only allows a single "_GoBack"
#    bookmark, so adding the new one removes the old one automatically,
#    leaving the final paragraph a plain empty paragraph.
# ---------------------------------------------------------------------------
$titleRng = $d.Content
$titleRng.Find.Execute("CONCLUSION GENERALE ET PERSPECTIVES", $true, $false, $false,
                        $false, $false, $true, 1, $false, "", 0)
$splitOffset = $titleRng.Start + "CONCLUSION GENERALE ET ".Length
$bmRng = $d.Range($splitOffset, $splitOffset)
$d.Bookmarks.Add("_GoBack", $bmRng)
